# Auto-generated edit script: updates Leve crafting-profit values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match refreshed
# market-board pricing data pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1149.7391
$ws.Range("J17").Value = 1161.5454
$ws.Range("L17").Value = 3484.6362
$ws.Range("N17").Value = -3820.6362
$ws.Range("H70").Value = 1345.64
$ws.Range("I70").Value = 1339.2084
$ws.Range("K70").Value = 4017.6252
$ws.Range("M70").Value = -3747.6252
$ws.Range("H73").Value = 1345.64
$ws.Range("I73").Value = 1339.2084
$ws.Range("K73").Value = 4017.6252
$ws.Range("M73").Value = -3081.6252
$ws.Range("H137").Value = 2175.5715
$ws.Range("I137").Value = 2839.125
$ws.Range("J137").Value = 1290.8334
$ws.Range("K137").Value = 8517.375
$ws.Range("L137").Value = 3872.5002
$ws.Range("M137").Value = -5967.375
$ws.Range("N137").Value = -8972.5002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 219.2
$ws.Range("I4").Value = 298.66666
$ws.Range("K4").Value = 298.66666
$ws.Range("M4").Value = -182.66666
$ws.Range("H5").Value = 108.666664
$ws.Range("I5").Value = 108.666664
$ws.Range("K5").Value = 108.666664
$ws.Range("M5").Value = 3.333336000000003
$ws.Range("H61").Value = 1904.9656
$ws.Range("I61").Value = 1840.8695
$ws.Range("J61").Value = 2150.6667
$ws.Range("K61").Value = 1840.8695
$ws.Range("L61").Value = 2150.6667
$ws.Range("M61").Value = -1628.8695
$ws.Range("N61").Value = -2574.6667
$ws.Range("H102").Value = 2119.64
$ws.Range("I102").Value = 1495.8823
$ws.Range("J102").Value = 3445.125
$ws.Range("K102").Value = 1495.8823
$ws.Range("L102").Value = 3445.125
$ws.Range("M102").Value = 126.1177
$ws.Range("N102").Value = -6689.125
$ws.Range("H132").Value = 5106.5
$ws.Range("I132").Value = 4121.227
$ws.Range("J132").Value = 7077.0454
$ws.Range("K132").Value = 12363.681
$ws.Range("L132").Value = 21231.1362
$ws.Range("M132").Value = -9833.681
$ws.Range("N132").Value = -26291.1362
$ws.Range("H136").Value = 1904.9656
$ws.Range("I136").Value = 1840.8695
$ws.Range("J136").Value = 2150.6667
$ws.Range("K136").Value = 5522.6085
$ws.Range("L136").Value = 6452.000100000001
$ws.Range("M136").Value = -2972.6085
$ws.Range("N136").Value = -11552.0001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 108.666664
$ws.Range("I4").Value = 108.666664
$ws.Range("K4").Value = 108.666664
$ws.Range("M4").Value = 6.333336000000003
$ws.Range("H22").Value = 820.6667
$ws.Range("I22").Value = 280.70587
$ws.Range("J22").Value = 10000
$ws.Range("K22").Value = 280.70587
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = -107.70587
$ws.Range("N22").Value = -10346
$ws.Range("H86").Value = 1701
$ws.Range("I86").Value = 1633.3334
$ws.Range("K86").Value = 1633.3334
$ws.Range("M86").Value = -510.3334
$ws.Range("H89").Value = 1701
$ws.Range("I89").Value = 1633.3334
$ws.Range("K89").Value = 8166.666999999999
$ws.Range("M89").Value = -2550.666999999999
$ws.Range("H99").Value = 1669.3667
$ws.Range("I99").Value = 1330
$ws.Range("J99").Value = 1895.6111
$ws.Range("K99").Value = 1330
$ws.Range("L99").Value = 1895.6111
$ws.Range("M99").Value = 168
$ws.Range("N99").Value = -4891.6111
$ws.Range("H105").Value = 2582.4285
$ws.Range("I105").Value = 2401.1765
$ws.Range("J105").Value = 3352.75
$ws.Range("K105").Value = 2401.1765
$ws.Range("L105").Value = 3352.75
$ws.Range("M105").Value = -654.1765
$ws.Range("N105").Value = -6846.75
$ws.Range("H134").Value = 2196.3713
$ws.Range("I134").Value = 1509.6279
$ws.Range("J134").Value = 3290.074
$ws.Range("K134").Value = 4528.8837
$ws.Range("L134").Value = 9870.222
$ws.Range("M134").Value = -1993.8837
$ws.Range("N134").Value = -14940.222

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 6775.4
$ws.Range("I7").Value = 11177.889
$ws.Range("J7").Value = 171.66667
$ws.Range("K7").Value = 11177.889
$ws.Range("L7").Value = 171.66667
$ws.Range("M7").Value = -11064.889
$ws.Range("N7").Value = -397.66667
$ws.Range("H22").Value = 610.2
$ws.Range("I22").Value = 300.5
$ws.Range("J22").Value = 816.6667
$ws.Range("K22").Value = 300.5
$ws.Range("L22").Value = 816.6667
$ws.Range("M22").Value = 49.5
$ws.Range("N22").Value = -1516.6667
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H62").Value = 13344.444
$ws.Range("I62").Value = 2528.5715
$ws.Range("J62").Value = 51200
$ws.Range("K62").Value = 2528.5715
$ws.Range("L62").Value = 51200
$ws.Range("M62").Value = -1904.5715
$ws.Range("N62").Value = -52448
$ws.Range("H65").Value = 13344.444
$ws.Range("I65").Value = 2528.5715
$ws.Range("J65").Value = 51200
$ws.Range("K65").Value = 12642.8575
$ws.Range("L65").Value = 256000
$ws.Range("M65").Value = -9522.8575
$ws.Range("N65").Value = -262240
$ws.Range("H132").Value = 1856.8462
$ws.Range("I132").Value = 1345.4783
$ws.Range("K132").Value = 4036.4349
$ws.Range("M132").Value = -1506.4349

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2633.476
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 2572.389
$ws.Range("K80").Value = 9000
$ws.Range("L80").Value = 7717.167
$ws.Range("M80").Value = -8064
$ws.Range("N80").Value = -9589.167000000001
$ws.Range("H83").Value = 2633.476
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 2572.389
$ws.Range("K83").Value = 27000
$ws.Range("L83").Value = 23151.501
$ws.Range("M83").Value = -22320
$ws.Range("N83").Value = -32511.501

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30450
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31560

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 20250
$ws.Range("J63").Value = 20250
$ws.Range("L63").Value = 20250
$ws.Range("N63").Value = -21498
$ws.Range("H66").Value = 20250
$ws.Range("J66").Value = 20250
$ws.Range("L66").Value = 60750
$ws.Range("N66").Value = -66990
